$wb = $excel.ActiveWorkbook
$active = $wb.ActiveSheet
$ws = $wb.Worksheets.Item("VTQaZ")

# Update row 6 (plugin hybrid vehicle), columns Q:AF (years 2035-2050) from 1 to 0
$ws.Range("Q6:AF6").Value = 0

# Update the selection on the VTQaZ sheet to match the saved selection state
$ws.Range("Q6:AF6").Select()

# Restore the originally active sheet/tab so the active tab selection is unchanged
$active.Activate()

$wb.Save()
